$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 196
$ws.Range("I2").Value = 196
$ws.Range("K2").Value = 196
$ws.Range("M2").Value = -83
$ws.Range("H15").Value = 470.125
$ws.Range("I15").Value = 470.125
$ws.Range("K15").Value = 1410.375
$ws.Range("M15").Value = -1241.375
$ws.Range("H70").Value = 2188.4736
$ws.Range("I70").Value = 2750.5
$ws.Range("J70").Value = 2038.6
$ws.Range("K70").Value = 8251.5
$ws.Range("L70").Value = 6115.799999999999
$ws.Range("M70").Value = -7981.5
$ws.Range("N70").Value = -6655.799999999999
$ws.Range("H73").Value = 2188.4736
$ws.Range("I73").Value = 2750.5
$ws.Range("J73").Value = 2038.6
$ws.Range("K73").Value = 8251.5
$ws.Range("L73").Value = 6115.799999999999
$ws.Range("M73").Value = -7315.5
$ws.Range("N73").Value = -7987.799999999999
$ws.Range("H94").Value = 1450
$ws.Range("I94").Value = 1450
$ws.Range("K94").Value = 1450
$ws.Range("M94").Value = -999
$ws.Range("H103").Value = 321.42856
$ws.Range("I103").Value = 321.42856
$ws.Range("K103").Value = 964.28568
$ws.Range("M103").Value = -378.28568
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H138").Value = 2036.9375
$ws.Range("I138").Value = 1287.8889
$ws.Range("K138").Value = 3863.6667
$ws.Range("M138").Value = 1276.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 33355
$ws.Range("J24").Value = 33355
$ws.Range("L24").Value = 33355
$ws.Range("N24").Value = -34103
$ws.Range("H45").Value = 11311.2
$ws.Range("I45").Value = 13705.75
$ws.Range("J45").Value = 1733
$ws.Range("K45").Value = 13705.75
$ws.Range("L45").Value = 1733
$ws.Range("M45").Value = -13328.75
$ws.Range("N45").Value = -2487
$ws.Range("H100").Value = 33355
$ws.Range("J100").Value = 33355
$ws.Range("L100").Value = 33355
$ws.Range("N100").Value = -35519
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 25370
$ws.Range("J105").Value = 25370
$ws.Range("L105").Value = 25370
$ws.Range("N105").Value = -32358
$ws.Range("H110").Value = 2221.1177
$ws.Range("I110").Value = 2113.3845
$ws.Range("J110").Value = 2571.25
$ws.Range("K110").Value = 2113.3845
$ws.Range("L110").Value = 2571.25
$ws.Range("M110").Value = -68.38450000000012
$ws.Range("N110").Value = -6661.25
$ws.Range("H132").Value = 2287.125
$ws.Range("I132").Value = 2256.7144
$ws.Range("K132").Value = 6770.1432
$ws.Range("M132").Value = -4240.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 925.7778
$ws.Range("J31").Value = 611.75
$ws.Range("L31").Value = 611.75
$ws.Range("N31").Value = -1201.75
$ws.Range("H34").Value = 925.7778
$ws.Range("J34").Value = 611.75
$ws.Range("L34").Value = 611.75
$ws.Range("N34").Value = -1015.75
$ws.Range("H105").Value = 627.53845
$ws.Range("I105").Value = 627.53845
$ws.Range("K105").Value = 627.53845
$ws.Range("M105").Value = 1119.46155
$ws.Range("H107").Value = 473
$ws.Range("I107").Value = 365
$ws.Range("J107").Value = 554
$ws.Range("K107").Value = 365
$ws.Range("L107").Value = 554
$ws.Range("M107").Value = 1555
$ws.Range("N107").Value = -4394
$ws.Range("H122").Value = 1202
$ws.Range("I122").Value = 1749.5
$ws.Range("J122").Value = 928.25
$ws.Range("K122").Value = 5248.5
$ws.Range("L122").Value = 2784.75
$ws.Range("M122").Value = -2798.5
$ws.Range("N122").Value = -7684.75
$ws.Range("H132").Value = 3749.8333
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560
$ws.Range("H134").Value = 3689.6296
$ws.Range("I134").Value = 2716.5833
$ws.Range("J134").Value = 4468.067
$ws.Range("K134").Value = 8149.749899999999
$ws.Range("L134").Value = 13404.201
$ws.Range("M134").Value = -5614.749899999999
$ws.Range("N134").Value = -18474.201

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 898.75
$ws.Range("I86").Value = 649.5
$ws.Range("J86").Value = 1148
$ws.Range("K86").Value = 1948.5
$ws.Range("L86").Value = 3444
$ws.Range("M86").Value = -762.5
$ws.Range("N86").Value = -5816
$ws.Range("H89").Value = 898.75
$ws.Range("I89").Value = 649.5
$ws.Range("J89").Value = 1148
$ws.Range("K89").Value = 5845.5
$ws.Range("L89").Value = 10332
$ws.Range("M89").Value = 82.5
$ws.Range("N89").Value = -22188
$ws.Range("H129").Value = 799.8570999999999
$ws.Range("I129").Value = 909.8
$ws.Range("J129").Value = 525
$ws.Range("K129").Value = 2729.4
$ws.Range("L129").Value = 1575
$ws.Range("M129").Value = 2270.6
$ws.Range("N129").Value = -11575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 750
$ws.Range("J4").Value = 750
$ws.Range("L4").Value = 750
$ws.Range("N4").Value = -974
$ws.Range("H80").Value = 2386.1428
$ws.Range("I80").Value = 2402
$ws.Range("K80").Value = 2402
$ws.Range("M80").Value = -1404
$ws.Range("H83").Value = 2386.1428
$ws.Range("I83").Value = 2402
$ws.Range("K83").Value = 12010
$ws.Range("M83").Value = -7018
$ws.Range("H126").Value = 5444.2856
$ws.Range("I126").Value = 5444.2856
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16332.8568
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -13862.8568
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2247.5
$ws.Range("I132").Value = 2211.4285
$ws.Range("K132").Value = 6634.2855
$ws.Range("M132").Value = -4104.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 40000
$ws.Range("I74").Value = 40000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 40000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -39002
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 40000
$ws.Range("I77").Value = 40000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 120000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -115008
$ws.Range("N77").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H136").Value = 2386.3333
$ws.Range("I136").Value = 2324
$ws.Range("K136").Value = 6972
$ws.Range("M136").Value = -4422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16999.166
$ws.Range("I132").Value = 24250
$ws.Range("K132").Value = 72750
$ws.Range("M132").Value = -70220
$ws.Range("H136").Value = 3630.6667
$ws.Range("I136").Value = 3547.4443
$ws.Range("J136").Value = 4130
$ws.Range("K136").Value = 10642.3329
$ws.Range("L136").Value = 12390
$ws.Range("M136").Value = -8092.332900000001
$ws.Range("N136").Value = -17490
